# Elimina EC anteriores y se agregan nuevos, se modifica base de datos
#
# The workers table (rows 16-28) is re-sorted: the block that belonged to
# "JIMENA CAROLINA BALLESTEROS BALLESTA" (CC 1235048451) moves down (to the
# bottom, rows 27-28) and the block for "DAIRO ALBERTO TRASLAVIÑA TORRES"
# (CC 1143397563) moves up (rows 17-26), with its period list reversed
# (2111 down to 2102). Row 16's "Salario Basico" also changes.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 16 - DAYANNA PATRICIA CALLE VILLARREAL keeps its place, only the
# Salario Basico (col G) changes.
$ws.Range("G16").Value = 877803

# Rows 17-26 - DAIRO ALBERTO TRASLAVIÑA TORRES, periods 2111 down to 2102
$ws.Range("C17").Value = "1143397563"
$ws.Range("D17").Value = "DAIRO ALBERTO TRASLAVIÑA TORRES"
$ws.Range("E17").Value = "2111"
$ws.Range("F17").Value = 26919
$ws.Range("G17").Value = 939249

$ws.Range("C18").Value = "1143397563"
$ws.Range("D18").Value = "DAIRO ALBERTO TRASLAVIÑA TORRES"
$ws.Range("E18").Value = "2110"
$ws.Range("F18").Value = 35112
$ws.Range("G18").Value = 939249

$ws.Range("E19").Value = "2109"
$ws.Range("E20").Value = "2108"
$ws.Range("E21").Value = "2107"
$ws.Range("E22").Value = "2106"
$ws.Range("E23").Value = "2105"
$ws.Range("E24").Value = "2104"
$ws.Range("E25").Value = "2103"
$ws.Range("E26").Value = "2102"

# Rows 27-28 - JIMENA CAROLINA BALLESTEROS BALLESTA, periods 1902 and 1901
$ws.Range("C27").Value = "1235048451"
$ws.Range("D27").Value = "JIMENA CAROLINA BALLESTEROS BALLESTA"
$ws.Range("E27").Value = "1902"
$ws.Range("F27").Value = 17667
$ws.Range("G27").Value = 877803

$ws.Range("C28").Value = "1235048451"
$ws.Range("D28").Value = "JIMENA CAROLINA BALLESTEROS BALLESTA"
$ws.Range("E28").Value = "1901"
$ws.Range("F28").Value = 1104
$ws.Range("G28").Value = 877803
